# Regenerate save_data: column G ("K") values were recomputed from the
# underlying strike (K) series instead of the old "Strike#" derivation,
# and std/mean were regenerated -> write the freshly calculated s_vals
# into column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 2
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 3
    26 = 1
    27 = 0
    28 = 2
    29 = 0
    30 = 1
    31 = 0
    32 = 3
    33 = 1
    34 = 0
    35 = 1
    36 = 0
    37 = 1
    38 = 0
    39 = 1
    40 = 2
    41 = 1
    42 = 0
    43 = 1
    44 = 2
    45 = 1
    46 = 2
    47 = 2
    48 = 2
    49 = 1
    50 = 1
    51 = 2
    52 = 1
    53 = 1
    54 = 1
    55 = 1
    56 = 1
    57 = 0
    58 = 1
    59 = 1
    60 = 0
    61 = 1
    63 = 0
    64 = 0
    65 = 0
    66 = 0
    67 = 0
    68 = 0
    69 = 1
    70 = 1
    71 = 2
    72 = 0
    73 = 1
    74 = 0
    75 = 1
    76 = 2
    77 = 0
    78 = 1
    80 = 1
    81 = 1
    84 = 1
    85 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
